$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $cell.Value2 + 1
}
